$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.7340267269605079
$ws.Cells.Item(2, 3).Value = 0.7680793040000253
$ws.Cells.Item(2, 4).Value = 0.5963715833518961
$ws.Cells.Item(2, 5).Value = 0.2075951406501559
$ws.Cells.Item(2, 7).Value = 0.002701280903774335
$ws.Cells.Item(2, 9).Value = 6.021838378453396
$ws.Cells.Item(2, 10).Value = 0.0828971149347808
$ws.Cells.Item(2, 11).Value = 1.445391888335763
$ws.Cells.Item(2, 13).Value = 0.6302291079695408

# Row 3
$ws.Cells.Item(3, 2).Value = 0.738542070033489
$ws.Cells.Item(3, 3).Value = 0.7430314300875409
$ws.Cells.Item(3, 4).Value = 0.5917118223406561
$ws.Cells.Item(3, 5).Value = 0.2061175895282794
$ws.Cells.Item(3, 7).Value = 0.002707340933266114
$ws.Cells.Item(3, 9).Value = 5.688095662163022
$ws.Cells.Item(3, 10).Value = 0.08238962656462689
$ws.Cells.Item(3, 11).Value = 1.430652675570428
$ws.Cells.Item(3, 13).Value = 0.6250162470212217

# Row 4
$ws.Cells.Item(4, 2).Value = 0.7423110122064003
$ws.Cells.Item(4, 3).Value = 0.7279255319842832
$ws.Cells.Item(4, 4).Value = 0.5891234622513224
$ws.Cells.Item(4, 5).Value = 0.2052959573583593
$ws.Cells.Item(4, 7).Value = 0.002711248263816709
$ws.Cells.Item(4, 9).Value = 5.482372429035081
$ws.Cells.Item(4, 10).Value = 0.08210526408340968
$ws.Cells.Item(4, 11).Value = 1.422894880972621
$ws.Cells.Item(4, 13).Value = 0.6222538927145891

# Row 5
$ws.Cells.Item(5, 2).Value = 0.7440963741009341
$ws.Cells.Item(5, 3).Value = 0.7218380456889122
$ws.Cells.Item(5, 4).Value = 0.588137156271145
$ws.Cells.Item(5, 5).Value = 0.2049826248164592
$ws.Cells.Item(5, 7).Value = 0.002712887599775284
$ws.Cells.Item(5, 9).Value = 5.398326946486975
$ws.Cells.Item(5, 10).Value = 0.08199622328440626
$ws.Cells.Item(5, 11).Value = 1.420057208277342
$ws.Cells.Item(5, 13).Value = 0.6212381459126064

# Row 6
$ws.Cells.Item(6, 2).Value = 0.7444078689928233
$ws.Cells.Item(6, 3).Value = 0.7208313331616125
$ws.Cells.Item(6, 4).Value = 0.5879775132323886
$ws.Cells.Item(6, 5).Value = 0.2049318931533612
$ws.Cells.Item(6, 7).Value = 0.002713162658438179
$ws.Cells.Item(6, 9).Value = 5.384358204187834
$ws.Cells.Item(6, 10).Value = 0.08197852993896504
$ws.Cells.Item(6, 11).Value = 1.419605529778494
$ws.Cells.Item(6, 13).Value = 0.6210761133040137

# Row 7
$ws.Cells.Item(7, 2).Value = 0.7423340814749793
$ws.Cells.Item(7, 3).Value = 0.7278431581735276
$ws.Cells.Item(7, 4).Value = 0.5891098834882769
$ws.Cells.Item(7, 5).Value = 0.2052916446787307
$ws.Cells.Item(7, 7).Value = 0.002711270181691287
$ws.Cells.Item(7, 9).Value = 5.481239831443446
$ws.Cells.Item(7, 10).Value = 0.08210376584169055
$ws.Cells.Item(7, 11).Value = 1.422855302100714
$ws.Cells.Item(7, 13).Value = 0.622239749254625

# Row 8
$ws.Cells.Item(8, 2).Value = 0.7353761265793821
$ws.Cells.Item(8, 3).Value = 0.7593856937097883
$ws.Cells.Item(8, 4).Value = 0.594708215051341
$ws.Cells.Item(8, 5).Value = 0.2070678954630907
$ws.Cells.Item(8, 7).Value = 0.002703331815378729
$ws.Cells.Item(8, 9).Value = 5.906924916223318
$ws.Cells.Item(8, 10).Value = 0.08271647344191635
$ws.Cells.Item(8, 11).Value = 1.440040865160853
$ws.Cells.Item(8, 13).Value = 0.628340547090211

# Row 9
$ws.Cells.Item(9, 2).Value = 0.7296900451270005
$ws.Cells.Item(9, 3).Value = 0.8234384501483873
$ws.Cells.Item(9, 4).Value = 0.6078573239166474
$ws.Cells.Item(9, 5).Value = 0.2112322006079843
$ws.Cells.Item(9, 7).Value = 0.002689235682881843
$ws.Cells.Item(9, 9).Value = 6.735780273134168
$ws.Cells.Item(9, 10).Value = 0.08413470533358947
$ws.Cells.Item(9, 11).Value = 1.484055622508606
$ws.Cells.Item(9, 13).Value = 0.643797794599152

# Row 10
$ws.Cells.Item(10, 2).Value = 0.73043921609937
$ws.Cells.Item(10, 3).Value = 0.8718813864604726
$ws.Cells.Item(10, 4).Value = 0.6188527073480543
$ws.Cells.Item(10, 5).Value = 0.2147103227314489
$ws.Cells.Item(10, 7).Value = 0.002679764182374342
$ws.Cells.Item(10, 9).Value = 7.341839690424422
$ws.Cells.Item(10, 10).Value = 0.08530984461799562
$ws.Cells.Item(10, 11).Value = 1.522774571969308
$ws.Cells.Item(10, 13).Value = 0.6573083566256699

# Row 11
$ws.Cells.Item(11, 2).Value = 0.7318666997935281
$ws.Cells.Item(11, 3).Value = 0.8942296735495461
$ws.Cells.Item(11, 4).Value = 0.6241473111240623
$ws.Cells.Item(11, 5).Value = 0.2163843165675701
$ws.Cells.Item(11, 7).Value = 0.002675644966615372
$ws.Cells.Item(11, 9).Value = 7.617079152355075
$ws.Cells.Item(11, 10).Value = 0.08587361157896112
$ws.Cells.Item(11, 11).Value = 1.541795836748008
$ws.Cells.Item(11, 13).Value = 0.6639280365216962

# Row 12
$ws.Cells.Item(12, 2).Value = 0.7325649024145378
$ws.Cells.Item(12, 3).Value = 0.9027378778812931
$ws.Cells.Item(12, 4).Value = 0.626194530705277
$ws.Cells.Item(12, 5).Value = 0.2170314698651339
$ws.Cells.Item(12, 7).Value = 0.002674112168704508
$ws.Cells.Item(12, 9).Value = 7.721250850931369
$ws.Cells.Item(12, 10).Value = 0.08609131053661656
$ws.Cells.Item(12, 11).Value = 1.549202810302575
$ws.Cells.Item(12, 13).Value = 0.6665032787689924

# Row 13
$ws.Cells.Item(13, 2).Value = 0.7324074982850277
$ws.Cells.Item(13, 3).Value = 0.9009034549124522
$ws.Cells.Item(13, 4).Value = 0.6257517426512607
$ws.Cells.Item(13, 5).Value = 0.2168915037791308
$ws.Cells.Item(13, 7).Value = 0.002674441083563895
$ws.Cells.Item(13, 9).Value = 7.698817909334707
$ws.Cells.Item(13, 10).Value = 0.08604423756408153
$ws.Cells.Item(13, 11).Value = 1.547598485755742
$ws.Cells.Item(13, 13).Value = 0.6659456005887137

# Row 14
$ws.Cells.Item(14, 2).Value = 0.7319209738560346
$ws.Cells.Item(14, 3).Value = 0.8949287352771194
$ws.Cells.Item(14, 4).Value = 0.6243148889877546
$ws.Cells.Item(14, 5).Value = 0.2164372924482549
$ws.Cells.Item(14, 7).Value = 0.002675518321045861
$ws.Cells.Item(14, 9).Value = 7.625650471648044
$ws.Cells.Item(14, 10).Value = 0.08589143728649162
$ws.Cells.Item(14, 11).Value = 1.542401114513098
$ws.Cells.Item(14, 13).Value = 0.6641385277957639

# Row 15
$ws.Cells.Item(15, 2).Value = 0.7316435352894928
$ws.Cells.Item(15, 3).Value = 0.8912749772208031
$ws.Cells.Item(15, 4).Value = 0.6234402843889768
$ws.Cells.Item(15, 5).Value = 0.2161608017337997
$ws.Cells.Item(15, 7).Value = 0.002676181678868961
$ws.Cells.Item(15, 9).Value = 7.580826392933659
$ws.Cells.Item(15, 10).Value = 0.08579839177569681
$ws.Cells.Item(15, 11).Value = 1.539244193717252
$ws.Cells.Item(15, 13).Value = 0.6630405781232938

# Row 16
$ws.Cells.Item(16, 2).Value = 0.7303679108658798
$ws.Cells.Item(16, 3).Value = 0.8704271912395711
$ws.Cells.Item(16, 4).Value = 0.618512597846518
$ws.Cells.Item(16, 5).Value = 0.2146027743826693
$ws.Cells.Item(16, 7).Value = 0.002680037179479709
$ws.Cells.Item(16, 9).Value = 7.323843792764791
$ws.Cells.Item(16, 10).Value = 0.08527358977968902
$ws.Cells.Item(16, 11).Value = 1.521559951445624
$ws.Cells.Item(16, 13).Value = 0.6568853102078549

# Row 17
$ws.Cells.Item(17, 2).Value = 0.72986465454008
$ws.Cells.Item(17, 3).Value = 0.8577179310487395
$ws.Cells.Item(17, 4).Value = 0.615564721080176
$ws.Cells.Item(17, 5).Value = 0.2136705169817787
$ws.Cells.Item(17, 7).Value = 0.002682450791563943
$ws.Cells.Item(17, 9).Value = 7.166082982046845
$ws.Cells.Item(17, 10).Value = 0.08495912760876934
$ws.Cells.Item(17, 11).Value = 1.511072935400421
$ws.Cells.Item(17, 13).Value = 0.6532308492943741

# Row 18
$ws.Cells.Item(18, 2).Value = 0.7296773938039394
$ws.Cells.Item(18, 3).Value = 0.8504371685233423
$ws.Cells.Item(18, 4).Value = 0.6138967245295817
$ws.Cells.Item(18, 5).Value = 0.2131429426607667
$ws.Cells.Item(18, 7).Value = 0.002683856877156393
$ws.Cells.Item(18, 9).Value = 7.075299191035128
$ws.Cells.Item(18, 10).Value = 0.08478100385874399
$ws.Cells.Item(18, 11).Value = 1.505173493994363
$ws.Cells.Item(18, 13).Value = 0.6511734580657205

# Row 19
$ws.Cells.Item(19, 2).Value = 0.7296315017385382
$ws.Cells.Item(19, 3).Value = 0.8479770355662595
$ws.Cells.Item(19, 4).Value = 0.6133366952095685
$ws.Cells.Item(19, 5).Value = 0.2129657967114795
$ws.Cells.Item(19, 7).Value = 0.002684336022954554
$ws.Cells.Item(19, 9).Value = 7.044553490592421
$ws.Cells.Item(19, 10).Value = 0.08472116549390307
$ws.Cells.Item(19, 11).Value = 1.503198738803519
$ws.Cells.Item(19, 13).Value = 0.6504845016081831

# Row 20
$ws.Cells.Item(20, 2).Value = 0.7299076399273758
$ws.Cells.Item(20, 3).Value = 0.8590678199424246
$ws.Cells.Item(20, 4).Value = 0.6158756756604475
$ws.Cells.Item(20, 5).Value = 0.213768863276087
$ws.Cells.Item(20, 7).Value = 0.002682192013588506
$ws.Cells.Item(20, 9).Value = 7.182881361375166
$ws.Cells.Item(20, 10).Value = 0.08499231827625664
$ws.Cells.Item(20, 11).Value = 1.512175580561689
$ws.Cells.Item(20, 13).Value = 0.6536152584445531

# Row 21
$ws.Cells.Item(21, 2).Value = 0.7320595882358703
$ws.Cells.Item(21, 3).Value = 0.8966824182757591
$ws.Cells.Item(21, 4).Value = 0.6247357791956745
$ws.Cells.Item(21, 5).Value = 0.216570345361049
$ws.Cells.Item(21, 7).Value = 0.002675201177196727
$ws.Cells.Item(21, 9).Value = 7.647142942594087
$ws.Cells.Item(21, 10).Value = 0.08593620397766699
$ws.Cells.Item(21, 11).Value = 1.543922158119557
$ws.Cells.Item(21, 13).Value = 0.6646674461818591

# Row 22
$ws.Cells.Item(22, 2).Value = 0.734385555408096
$ws.Cells.Item(22, 3).Value = 0.9215306185523673
$ws.Cells.Item(22, 4).Value = 0.6307728244549651
$ws.Cells.Item(22, 5).Value = 0.2184785213946796
$ws.Cells.Item(22, 7).Value = 0.002670789905875342
$ws.Cells.Item(22, 9).Value = 7.950246226551144
$ws.Cells.Item(22, 10).Value = 0.08657764989490602
$ws.Cells.Item(22, 11).Value = 1.565860415766394
$ws.Cells.Item(22, 13).Value = 0.6722902330855121

# Row 23
$ws.Cells.Item(23, 2).Value = 0.7330595390218946
$ws.Cells.Item(23, 3).Value = 0.9082442196991565
$ws.Cells.Item(23, 4).Value = 0.6275281305436522
$ws.Cells.Item(23, 5).Value = 0.2174530066632414
$ws.Cells.Item(23, 7).Value = 0.002673129918552736
$ws.Cells.Item(23, 9).Value = 7.788499878997527
$ws.Cells.Item(23, 10).Value = 0.08623304552454414
$ws.Cells.Item(23, 11).Value = 1.554042128986993
$ws.Cells.Item(23, 13).Value = 0.6681851136642933

# Row 24
$ws.Cells.Item(24, 2).Value = 0.7298878884186308
$ws.Cells.Item(24, 3).Value = 0.8584574543587564
$ws.Cells.Item(24, 4).Value = 0.6157350097032008
$ws.Cells.Item(24, 5).Value = 0.213724374778316
$ws.Cells.Item(24, 7).Value = 0.002682308949738761
$ws.Cells.Item(24, 9).Value = 7.175287079636632
$ws.Cells.Item(24, 10).Value = 0.08497730447680851
$ws.Cells.Item(24, 11).Value = 1.511676670867899
$ws.Cells.Item(24, 13).Value = 0.6534413312717575

# Row 25
$ws.Cells.Item(25, 2).Value = 0.7303686291344889
$ws.Cells.Item(25, 3).Value = 0.8058707437891712
$ws.Cells.Item(25, 4).Value = 0.6040665947994341
$ws.Cells.Item(25, 5).Value = 0.2100323874875087
$ws.Cells.Item(25, 7).Value = 0.002692892804611373
$ws.Cells.Item(25, 9).Value = 6.512113746668092
$ws.Cells.Item(25, 10).Value = 0.08372772918013638
$ws.Cells.Item(25, 11).Value = 1.471034760565544
$ws.Cells.Item(25, 13).Value = 0.6392397168833526

Write-Output "Updated 216 cells"